$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Simple price-only updates (column D holds numeric-looking text) ---
Set-TextValue "D2"  "231.77"
Set-TextValue "D3"  "22.85"
Set-TextValue "D4"  "5.583"
Set-TextValue "D5"  "0.05588"
Set-TextValue "D6"  "3.427"
Set-TextValue "D7"  "6.501"
Set-TextValue "D8"  "1.167"
Set-TextValue "D9"  "0.7981"
Set-TextValue "D10" "0.1416"
Set-TextValue "D11" "0.07415"
Set-TextValue "D12" "0.03153"
Set-TextValue "D13" "0.02949"
Set-TextValue "D14" "0.09264"
Set-TextValue "D15" "0.001655"
Set-TextValue "D16" "3.290"
Set-TextValue "D17" "0.04735"

# --- Row 18: TigerCash -> One ---
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005943"
$ws.Range("E18").Value = "17OneONE"

# --- Row 19: HotbitToken -> TigerCash ---
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006264"
$ws.Range("E19").Value = "18TigerCashTCH"

# --- Row 20: BitKan -> HotbitToken ---
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005270"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# --- Row 21: NitroEx -> BitKan ---
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001071"
$ws.Range("E21").Value = "20BitKanKAN"

# --- Row 22: LEO -> NitroEx ---
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001512"
$ws.Range("E22").Value = "21NitroExNTX"

# --- Row 23: BTSEToken -> LEO ---
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.680"
$ws.Range("E23").Value = "22LEOLEO"

# --- Row 24: One -> BTSEToken ---
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.192"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- Remaining simple price-only updates ---
Set-TextValue "D25" "0.3322"

Set-TextValue "D27" "0.0008375"

Set-TextValue "D40" "0.04085"
Set-TextValue "D41" "0.007185"
Set-TextValue "D42" "0.003404"
Set-TextValue "D43" "0.1034"

Set-TextValue "D44" "0.008195"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"

$ws.Range("E45").Value = "44ACDXExchangeACXT"

Set-TextValue "D46" "0.00005575"
Set-TextValue "D47" "0.00000000756"
Set-TextValue "D48" "0.6805"
Set-TextValue "D49" "0.09537"
Set-TextValue "D50" "0.00002116"
Set-TextValue "D51" "0.01018"
